$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Altman's Z
$ws.Range("B2").Value = 0.3507
$ws.Range("C2").Value = "'0.2435"
$ws.Range("D2").Value = "'0.3507"
$ws.Range("E2").Value = 0.2751
$ws.Range("F2").Value = 0.8022

# Row 3 - Financial Variables and Sector
$ws.Range("B3").Value = 0.6484
$ws.Range("C3").Value = "'0.6464"
$ws.Range("D3").Value = "'0.6484"
$ws.Range("E3").Value = 0.6457000000000001
$ws.Range("F3").Value = 0.9432

# Row 4 - Financial Variables, Sector, and NLP Features
$ws.Range("B4").Value = 0.6556999999999999
$ws.Range("C4").Value = "'0.6573"
$ws.Range("D4").Value = "'0.6557"
$ws.Range("E4").Value = 0.6536999999999999
$ws.Range("F4").Value = 0.946

# Row 5 - Majority Baseline
$ws.Range("B5").Value = 0.3013
